# Update "想去人数" (number of people interested) values in the
# "展览" and "全部类型" worksheets to match the newly scraped data.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes    = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates (column F)
$sheetExhibition.Range("F2").Value = 5576
$sheetExhibition.Range("F3").Value = 12978
$sheetExhibition.Range("F5").Value = 634
$sheetExhibition.Range("F7").Value = 405
$sheetExhibition.Range("F8").Value = 1227

# 全部类型 sheet updates (column F)
$sheetAllTypes.Range("F2").Value  = 5576
$sheetAllTypes.Range("F4").Value  = 12978
$sheetAllTypes.Range("F6").Value  = 634
$sheetAllTypes.Range("F10").Value = 405
$sheetAllTypes.Range("F11").Value = 1227
